# Scheduled runner update: refresh cached market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns) for the rows
# whose underlying price data changed, across the relevant sheets.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1355.75
$ws.Range("I38").Value = 369.6
$ws.Range("J38").Value = 2999.3333
$ws.Range("K38").Value = 1108.8
$ws.Range("L38").Value = 8997.999899999999
$ws.Range("M38").Value = -736.8000000000002
$ws.Range("N38").Value = -9741.999899999999
$ws.Range("H53").Value = 153.12
$ws.Range("I53").Value = 105.63636
$ws.Range("J53").Value = 190.42857
$ws.Range("K53").Value = 105.63636
$ws.Range("L53").Value = 190.42857
$ws.Range("M53").Value = 531.36364
$ws.Range("N53").Value = -1464.42857
$ws.Range("H113").Value = 2701.7646
$ws.Range("I113").Value = 2111
$ws.Range("J113").Value = 2947.9167
$ws.Range("K113").Value = 2111
$ws.Range("L113").Value = 2947.9167
$ws.Range("M113").Value = 1143
$ws.Range("N113").Value = -9455.9167
$ws.Range("H116").Value = 1914.9667
$ws.Range("I116").Value = 1748.2
$ws.Range("K116").Value = 1748.2
$ws.Range("M116").Value = 1693.8
$ws.Range("H132").Value = 2761.2
$ws.Range("I132").Value = 2929.111
$ws.Range("J132").Value = 1250
$ws.Range("K132").Value = 8787.332999999999
$ws.Range("L132").Value = 3750
$ws.Range("M132").Value = -6257.332999999999
$ws.Range("N132").Value = -8810
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -70120
$ws.Range("H137").Value = 1571.6852
$ws.Range("I137").Value = 1258.0238
$ws.Range("J137").Value = 2669.5
$ws.Range("K137").Value = 3774.0714
$ws.Range("L137").Value = 8008.5
$ws.Range("M137").Value = -1224.0714
$ws.Range("N137").Value = -13108.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8261.825999999999
$ws.Range("I61").Value = 5375.0415
$ws.Range("J61").Value = 11411.046
$ws.Range("K61").Value = 5375.0415
$ws.Range("L61").Value = 11411.046
$ws.Range("M61").Value = -5163.0415
$ws.Range("N61").Value = -11835.046
$ws.Range("H132").Value = 2064.1912
$ws.Range("I132").Value = 1606.9246
$ws.Range("K132").Value = 4820.7738
$ws.Range("M132").Value = -2290.7738
$ws.Range("H136").Value = 8261.825999999999
$ws.Range("I136").Value = 5375.0415
$ws.Range("J136").Value = 11411.046
$ws.Range("K136").Value = 16125.1245
$ws.Range("L136").Value = 34233.138
$ws.Range("M136").Value = -13575.1245
$ws.Range("N136").Value = -39333.138
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 35649.633
$ws.Range("I134").Value = 2187.2856
$ws.Range("K134").Value = 6561.8568
$ws.Range("M134").Value = -4026.8568
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3573.9333
$ws.Range("I62").Value = 3280.4546
$ws.Range("J62").Value = 4381
$ws.Range("K62").Value = 3280.4546
$ws.Range("L62").Value = 4381
$ws.Range("M62").Value = -2656.4546
$ws.Range("N62").Value = -5629
$ws.Range("H65").Value = 3573.9333
$ws.Range("I65").Value = 3280.4546
$ws.Range("J65").Value = 4381
$ws.Range("K65").Value = 16402.273
$ws.Range("L65").Value = 21905
$ws.Range("M65").Value = -13282.273
$ws.Range("N65").Value = -28145
$ws.Range("H68").Value = 39628.332
$ws.Range("J68").Value = 39628.332
$ws.Range("L68").Value = 39628.332
$ws.Range("N68").Value = -41126.332
$ws.Range("H71").Value = 39628.332
$ws.Range("J71").Value = 39628.332
$ws.Range("L71").Value = 118884.996
$ws.Range("N71").Value = -126372.996
$ws.Range("H99").Value = 3779.0833
$ws.Range("I99").Value = 3304.5
$ws.Range("J99").Value = 4728.25
$ws.Range("K99").Value = 3304.5
$ws.Range("L99").Value = 4728.25
$ws.Range("M99").Value = -1806.5
$ws.Range("N99").Value = -7724.25
$ws.Range("H126").Value = 3779.0833
$ws.Range("I126").Value = 3304.5
$ws.Range("J126").Value = 4728.25
$ws.Range("K126").Value = 9913.5
$ws.Range("L126").Value = 14184.75
$ws.Range("M126").Value = -7443.5
$ws.Range("N126").Value = -19124.75
$ws.Range("H132").Value = 3939.9807
$ws.Range("I132").Value = 4738.9116
$ws.Range("K132").Value = 14216.7348
$ws.Range("M132").Value = -11686.7348
$ws.Range("H134").Value = 3060.8667
$ws.Range("I134").Value = 2461.7222
$ws.Range("J134").Value = 3460.2964
$ws.Range("K134").Value = 7385.1666
$ws.Range("L134").Value = 10380.8892
$ws.Range("M134").Value = -4850.1666
$ws.Range("N134").Value = -15450.8892
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 504.45456
$ws.Range("I7").Value = 558.3333
$ws.Range("J7").Value = 439.8
$ws.Range("K7").Value = 1674.9999
$ws.Range("L7").Value = 1319.4
$ws.Range("M7").Value = -1562.9999
$ws.Range("N7").Value = -1543.4
$ws.Range("H41").Value = 79.5
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H113").Value = 699.2643399999999
$ws.Range("I113").Value = 695.7183
$ws.Range("J113").Value = 715
$ws.Range("K113").Value = 2087.1549
$ws.Range("L113").Value = 2145
$ws.Range("M113").Value = 82.8451
$ws.Range("N113").Value = -6485
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3332.5
$ws.Range("I80").Value = 2750
$ws.Range("J80").Value = 3449
$ws.Range("K80").Value = 2750
$ws.Range("L80").Value = 3449
$ws.Range("M80").Value = -1752
$ws.Range("N80").Value = -5445
$ws.Range("H83").Value = 3332.5
$ws.Range("I83").Value = 2750
$ws.Range("J83").Value = 3449
$ws.Range("K83").Value = 13750
$ws.Range("L83").Value = 17245
$ws.Range("M83").Value = -8758
$ws.Range("N83").Value = -27229
$ws.Range("H113").Value = 2683.8948
$ws.Range("I113").Value = 3428.3333
$ws.Range("K113").Value = 3428.3333
$ws.Range("M113").Value = -1258.3333
$ws.Range("H132").Value = 5535.3794
$ws.Range("I132").Value = 2247.92
$ws.Range("K132").Value = 6743.76
$ws.Range("M132").Value = -4213.76
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6562.9204
$ws.Range("I122").Value = 6131.5684
$ws.Range("J122").Value = 7561.8423
$ws.Range("K122").Value = 18394.7052
$ws.Range("L122").Value = 22685.5269
$ws.Range("M122").Value = -15944.7052
$ws.Range("N122").Value = -27585.5269
$ws.Range("H136").Value = 3496.4707
$ws.Range("I136").Value = 1895.5555
$ws.Range("J136").Value = 6628.696
$ws.Range("K136").Value = 5686.666499999999
$ws.Range("L136").Value = 19886.088
$ws.Range("M136").Value = -3136.666499999999
$ws.Range("N136").Value = -24986.088
